# Update the "債務" (Debt) worksheet (sheet index 6) to match the standard
# property-record layout used by the other sheets: add header labels in row 1,
# and append legislator_name / legislator_id / source_file / index columns
# (H:N) to the data row, relabeling the existing B:G header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# --- Row 1: header labels ---
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: data values ---
$ws.Range("B2").Value = "抵押"
$ws.Range("C2").Value = "李桐豪"
$ws.Range("D2").Value = "合作金庫銀行西門分行臺北市萬華區昆明街"
$ws.Range("E2").Value = 313034
$ws.Range("F2").Value = "83年07月27日"
$ws.Range("G2").Value = "公教貸款"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-25"
$ws.Range("K2").Value = "李桐豪"
$ws.Range("L2").Value = 896
$ws.Range("M2").Value = "tmpe99a1"
$ws.Range("N2").Value = 99

# --- Apply header (bold/bordered) style to the new header cells, and the
#     plain data style to the new data cells, matching the existing columns ---
$ws.Range("H1:N1").Style = $ws.Range("B1").Style
$ws.Range("H2:N2").Style = $ws.Range("B2").Style

# --- Also update "具有相當價值之財產" (sheet 5) which shares the renamed
#     string "otherbonds" -> "antique" in the shared-strings table; its
#     property_category cell (F2) keeps pointing at that same string, which
#     now reads "antique" automatically. No direct write needed there.
